$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Fix the "Bursa" column (E): values were entered x1000 too large
#    (e.g. 2000000 instead of 2000 lei) - scale them back down.
#    Row 14 is a genuine data correction (was mis-keyed) -> 1000.
# ------------------------------------------------------------------
$rows = @(4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,64,65,66,67,68,69,70,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96)
$vals = @(2000,1500,2000,1500,500,1500,1500,1500,1500,1500,1000,1500,500,500,1500,500,500,1500,1500,1500,2500,1500,1500,500,1500,1500,1500,1500,1500,1500,1500,1500,2500,1500,2000,1500,1500,1500,1500,1500,1500,1500,500,2000,2000,1500,500,1500,2000,1500,1500,1500,500,2000,500,1500,1500,2500,2500,1500,1500,2500,1500,500,1500,1500,1500,1500,1500,1500,2000,1500,1500,1500,500,1500,1500,500,1500,2000,1500,1500,1500,500,1500,1500,2500,2000,1500,1500)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 5).Value = $vals[$i]
}

# ------------------------------------------------------------------
# 2. Add the new instruction text under the other "Folosind functiile..."
#    helper notes, in a merged, highlighted cell G17:I17.
# ------------------------------------------------------------------
$noteRange = $ws.Range("G17:I17")
$noteRange.Merge() | Out-Null
$noteRange.Font.Bold = $true
$noteRange.Font.Size = 10
$noteRange.Font.Color = 2110976
$g17 = $ws.Range("G17")
$g17.Value = "Sortati datele descrescator in functie de media anuala"
$g17.Interior.PatternColorIndex = -4105
$g17.Interior.ThemeColor = 0
$g17.Interior.TintAndShade = 0
$g17.Interior.Pattern = 1

# ------------------------------------------------------------------
# 3. Emphasise the "Bursa cea mai mica" note (G13) a bit more.
# ------------------------------------------------------------------
$g13 = $ws.Range("G13")
$g13.Font.Bold = $true
$g13.Font.Size = 11
$g13.Font.ColorIndex = 18
$ws.Rows.Item(13).RowHeight = 15

# ------------------------------------------------------------------
# 4. Move the active selection to the newly added note.
# ------------------------------------------------------------------
$ws.Range("G17:I17").Select() | Out-Null
